$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 3 (pushes the existing row 3 and everything
# below it down by one), then populate it with the new day-2 (June/2025)
# faturamento entry: Dia=2, total_venda=27652.8, Mes=6, Ano=2025, Periodo="06/2025"
$ws.Rows.Item(3).Insert()

$ws.Cells.Item(3, 1).Value = 2
$ws.Cells.Item(3, 2).Value = 27652.8
$ws.Cells.Item(3, 3).Value = 6
$ws.Cells.Item(3, 4).Value = 2025
$ws.Cells.Item(3, 5).Value = "06/2025"
